# Update "想去人数" (F column) figures that were refreshed in the scraped data.
# Source: gh-pages data regeneration commit (456a3b4).

$wb = $excel.ActiveWorkbook

# ---- Sheet: 展览 (sheetId 1) ----
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 9756
$ws1.Range("F4").Value  = 2517
$ws1.Range("F8").Value  = 477
$ws1.Range("F9").Value  = 714
$ws1.Range("F11").Value = 1214
$ws1.Range("F12").Value = 1011
$ws1.Range("F13").Value = 3016
$ws1.Range("F14").Value = 2286
$ws1.Range("F16").Value = 1985
$ws1.Range("F21").Value = 317
$ws1.Range("F22").Value = 26
$ws1.Range("F23").Value = 199
$ws1.Range("F25").Value = 37
$ws1.Range("F28").Value = 326
$ws1.Range("F29").Value = 537
$ws1.Range("F30").Value = 39
$ws1.Range("F31").Value = 176
$ws1.Range("F32").Value = 1553
$ws1.Range("F33").Value = 215
$ws1.Range("F34").Value = 1550
$ws1.Range("F35").Value = 69
$ws1.Range("F36").Value = 370
$ws1.Range("F37").Value = 36
$ws1.Range("F38").Value = 401
$ws1.Range("F39").Value = 818
$ws1.Range("F40").Value = 79
$ws1.Range("F41").Value = 323

# ---- Sheet: 演出 (sheetId 2) ----
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F4").Value = 30

# ---- Sheet: 全部类型 (sheetId 4) ----
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 9756
$ws4.Range("F4").Value  = 2517
$ws4.Range("F10").Value = 477
$ws4.Range("F11").Value = 714
$ws4.Range("F13").Value = 1214
$ws4.Range("F14").Value = 1011
$ws4.Range("F15").Value = 3016
$ws4.Range("F16").Value = 2286
$ws4.Range("F18").Value = 1985
$ws4.Range("F23").Value = 317
$ws4.Range("F24").Value = 26
$ws4.Range("F25").Value = 199
$ws4.Range("F27").Value = 37
$ws4.Range("F30").Value = 326
$ws4.Range("F31").Value = 537
$ws4.Range("F32").Value = 30
$ws4.Range("F35").Value = 39
$ws4.Range("F36").Value = 176
$ws4.Range("F37").Value = 1553
$ws4.Range("F39").Value = 216
$ws4.Range("F40").Value = 1550
$ws4.Range("F41").Value = 69
$ws4.Range("F43").Value = 370
$ws4.Range("F44").Value = 36
$ws4.Range("F45").Value = 401
$ws4.Range("F46").Value = 818
$ws4.Range("F47").Value = 79
$ws4.Range("F48").Value = 323
